$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1855203619909502
$ws.Range("C2").Value = 0.6018099547511312
$ws.Range("J2").Value = 0.01357466063348416
$ws.Range("P2").Value = 0.1266968325791855
$ws.Range("S2").Value = 0.07239819004524888
$ws.Range("J3").Value = 0.01503759398496241
$ws.Range("P3").Value = 0.7894736842105263
$ws.Range("S3").Value = 0.1954887218045113
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.59375
$ws.Range("S4").Value = 0.34375
$ws.Range("B6").Value = 0.03982300884955752
$ws.Range("D6").Value = 0.02212389380530973
$ws.Range("F6").Value = 0.08849557522123894
$ws.Range("J6").Value = 0.2389380530973451
$ws.Range("O6").Value = 0.01769911504424779
$ws.Range("Q6").Value = 0.1814159292035398
$ws.Range("R6").Value = 0.0752212389380531
$ws.Range("S6").Value = 0.336283185840708
$ws.Range("B7").Value = 0.08808290155440414
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.08808290155440414
$ws.Range("J7").Value = 0.1295336787564767
$ws.Range("O7").Value = 0.05181347150259067
$ws.Range("Q7").Value = 0.1243523316062176
$ws.Range("R7").Value = 0.1295336787564767
$ws.Range("S7").Value = 0.383419689119171
$ws.Range("B8").Value = 0.07163323782234957
$ws.Range("D8").Value = 0.02292263610315186
$ws.Range("E8").Value = 0.002865329512893983
$ws.Range("F8").Value = 0.06876790830945559
$ws.Range("J8").Value = 0.1375358166189112
$ws.Range("O8").Value = 0.02292263610315186
$ws.Range("Q8").Value = 0.1547277936962751
$ws.Range("R8").Value = 0.1432664756446991
$ws.Range("S8").Value = 0.3753581661891118
$ws.Range("B9").Value = 0.08695652173913043
$ws.Range("D9").Value = 0.00966183574879227
$ws.Range("F9").Value = 0.08695652173913043
$ws.Range("J9").Value = 0.1304347826086956
$ws.Range("O9").Value = 0.02415458937198068
$ws.Range("Q9").Value = 0.1497584541062802
$ws.Range("R9").Value = 0.1449275362318841
$ws.Range("S9").Value = 0.3671497584541063
$ws.Range("B10").Value = 0.08684863523573201
$ws.Range("D10").Value = 0.01488833746898263
$ws.Range("E10").Value = 0.0008271298593879239
$ws.Range("F10").Value = 0.07775020678246485
$ws.Range("J10").Value = 0.1066997518610422
$ws.Range("O10").Value = 0.02233250620347394
$ws.Range("Q10").Value = 0.2150537634408602
$ws.Range("R10").Value = 0.1133167907361456
$ws.Range("S10").Value = 0.3622828784119106
$ws.Range("G11").Value = 0.1353135313531353
$ws.Range("J11").Value = 0.06930693069306931
$ws.Range("K11").Value = 0.1947194719471947
$ws.Range("L11").Value = 0.5874587458745875
$ws.Range("S11").Value = 0.0132013201320132
$ws.Range("G12").Value = 0.732620320855615
$ws.Range("J12").Value = 0.1978609625668449
$ws.Range("K12").Value = 0.0053475935828877
$ws.Range("L12").Value = 0.0374331550802139
$ws.Range("S12").Value = 0.0267379679144385
$ws.Range("G13").Value = 0.6774193548387096
$ws.Range("J13").Value = 0.2580645161290323
$ws.Range("S13").Value = 0.06451612903225806
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.04035874439461883
$ws.Range("H15").Value = 0.1479820627802691
$ws.Range("I15").Value = 0.05829596412556054
$ws.Range("J15").Value = 0.2869955156950673
$ws.Range("K15").Value = 0.07174887892376682
$ws.Range("M15").Value = 0.004484304932735426
$ws.Range("O15").Value = 0.08071748878923767
$ws.Range("S15").Value = 0.3094170403587444
$ws.Range("F16").Value = 0.01342281879194631
$ws.Range("H16").Value = 0.1543624161073825
$ws.Range("I16").Value = 0.04697986577181208
$ws.Range("J16").Value = 0.436241610738255
$ws.Range("K16").Value = 0.1073825503355705
$ws.Range("M16").Value = 0.006711409395973154
$ws.Range("O16").Value = 0.02684563758389262
$ws.Range("S16").Value = 0.2080536912751678
$ws.Range("F17").Value = 0.009828009828009828
$ws.Range("H17").Value = 0.2113022113022113
$ws.Range("I17").Value = 0.09582309582309582
$ws.Range("J17").Value = 0.4127764127764127
$ws.Range("K17").Value = 0.09828009828009827
$ws.Range("M17").Value = 0.009828009828009828
$ws.Range("N17").Value = 0.002457002457002457
$ws.Range("O17").Value = 0.04176904176904177
$ws.Range("S17").Value = 0.1179361179361179
$ws.Range("F18").Value = 0.01171875
$ws.Range("H18").Value = 0.16015625
$ws.Range("I18").Value = 0.109375
$ws.Range("J18").Value = 0.42578125
$ws.Range("K18").Value = 0.11328125
$ws.Range("M18").Value = 0.01953125
$ws.Range("O18").Value = 0.05078125
$ws.Range("S18").Value = 0.109375
$ws.Range("F19").Value = 0.0141718334809566
$ws.Range("H19").Value = 0.1514614703277236
$ws.Range("I19").Value = 0.1080602302922941
$ws.Range("J19").Value = 0.4047829937998229
$ws.Range("K19").Value = 0.1231178033658104
$ws.Range("M19").Value = 0.01860053144375554
$ws.Range("O19").Value = 0.07971656333038087
$ws.Range("S19").Value = 0.100088573959256
